# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values for rows 2-30 (replacing the prior Strike# values).
$newK = @{
    2  = 3
    3  = 4
    4  = 6
    5  = 7
    6  = 5
    7  = 3
    8  = 11
    9  = 4
    10 = 5
    11 = 10
    12 = 8
    13 = 9
    14 = 8
    15 = 8
    16 = 9
    17 = 8
    18 = 8
    19 = 10
    20 = 13
    21 = 5
    22 = 8
    23 = 12
    24 = 8
    25 = 7
    26 = 9
    27 = 6
    28 = 5
    29 = 4
    30 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
